# Auto-generated from diff: apply cell value updates across 8 sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1800
$ws.Range("I18").Value = 1800
$ws.Range("K18").Value = 1800
$ws.Range("M18").Value = -1516
$ws.Range("H43").Value = 3473.7273
$ws.Range("I43").Value = 2698.5
$ws.Range("K43").Value = 2698.5
$ws.Range("M43").Value = -2629.5
$ws.Range("H48").Value = 5000
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H56").Value = 5000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H74").Value = 5605.6313
$ws.Range("I74").Value = 5605.6313
$ws.Range("K74").Value = 5605.6313
$ws.Range("M74").Value = -4669.6313
$ws.Range("H76").Value = 5049.9
$ws.Range("I76").Value = 4937.375
$ws.Range("K76").Value = 4937.375
$ws.Range("M76").Value = -4622.375
$ws.Range("H77").Value = 5605.6313
$ws.Range("I77").Value = 5605.6313
$ws.Range("K77").Value = 28028.1565
$ws.Range("M77").Value = -23348.1565
$ws.Range("H79").Value = 5049.9
$ws.Range("I79").Value = 4937.375
$ws.Range("K79").Value = 4937.375
$ws.Range("M79").Value = -3845.375
$ws.Range("H116").Value = 102479.8
$ws.Range("I116").Value = 252050.75
$ws.Range("J116").Value = 2765.8333
$ws.Range("K116").Value = 252050.75
$ws.Range("L116").Value = 2765.8333
$ws.Range("M116").Value = -248608.75
$ws.Range("N116").Value = -9649.8333
$ws.Range("H132").Value = 53973.79
$ws.Range("I132").Value = 56805.332
$ws.Range("K132").Value = 170415.996
$ws.Range("M132").Value = -167885.996

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4057.1082
$ws.Range("I32").Value = 2503.1667
$ws.Range("K32").Value = 2503.1667
$ws.Range("M32").Value = -2216.1667
$ws.Range("H97").Value = 981.8095
$ws.Range("I97").Value = 624.6
$ws.Range("K97").Value = 624.6
$ws.Range("M97").Value = -128.6
$ws.Range("H132").Value = 14709088
$ws.Range("I132").Value = 1965.0385
$ws.Range("K132").Value = 5895.1155
$ws.Range("M132").Value = -3365.1155

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2627.8
$ws.Range("I105").Value = 2175.2856
$ws.Range("J105").Value = 3203.7273
$ws.Range("K105").Value = 2175.2856
$ws.Range("L105").Value = 3203.7273
$ws.Range("M105").Value = -428.2856000000002
$ws.Range("N105").Value = -6697.7273

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1240.2307
$ws.Range("I16").Value = 929.36365
$ws.Range("K16").Value = 929.36365
$ws.Range("M16").Value = -642.36365
$ws.Range("H31").Value = 2376.3809
$ws.Range("I31").Value = 2299.3572
$ws.Range("J31").Value = 2530.4285
$ws.Range("K31").Value = 2299.3572
$ws.Range("L31").Value = 2530.4285
$ws.Range("M31").Value = -2004.3572
$ws.Range("N31").Value = -3120.4285
$ws.Range("H34").Value = 2376.3809
$ws.Range("I34").Value = 2299.3572
$ws.Range("J34").Value = 2530.4285
$ws.Range("K34").Value = 2299.3572
$ws.Range("L34").Value = 2530.4285
$ws.Range("M34").Value = -2097.3572
$ws.Range("N34").Value = -2934.4285
$ws.Range("H58").Value = 1461.381
$ws.Range("I58").Value = 1551
$ws.Range("K58").Value = 1551
$ws.Range("M58").Value = -1348
$ws.Range("H105").Value = 1942.5555
$ws.Range("I105").Value = 1694.0714
$ws.Range("K105").Value = 1694.0714
$ws.Range("M105").Value = 52.92859999999996
$ws.Range("H107").Value = 639.14813
$ws.Range("I107").Value = 620.3684
$ws.Range("J107").Value = 683.75
$ws.Range("K107").Value = 620.3684
$ws.Range("L107").Value = 683.75
$ws.Range("M107").Value = 1299.6316
$ws.Range("N107").Value = -4523.75
$ws.Range("H113").Value = 1240.2307
$ws.Range("I113").Value = 929.36365
$ws.Range("K113").Value = 929.36365
$ws.Range("M113").Value = 1240.63635
$ws.Range("H136").Value = 1461.381
$ws.Range("I136").Value = 1551
$ws.Range("K136").Value = 4653
$ws.Range("M136").Value = -2103

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 958.75
$ws.Range("I8").Value = 958.75
$ws.Range("K8").Value = 2876.25
$ws.Range("M8").Value = -2737.25
$ws.Range("H25").Value = 14989
$ws.Range("J25").Value = 14989
$ws.Range("L25").Value = 44967
$ws.Range("N25").Value = -45305
$ws.Range("H30").Value = 14989
$ws.Range("J30").Value = 14989
$ws.Range("L30").Value = 44967
$ws.Range("N30").Value = -45171
$ws.Range("H51").Value = 1341.4
$ws.Range("I51").Value = 1303
$ws.Range("K51").Value = 3909
$ws.Range("M51").Value = -3449
$ws.Range("H56").Value = 8620.154
$ws.Range("I56").Value = 8620.154
$ws.Range("K56").Value = 8620.154
$ws.Range("M56").Value = -8090.154
$ws.Range("H64").Value = 1011.5
$ws.Range("I64").Value = 1011.5
$ws.Range("K64").Value = 3034.5
$ws.Range("M64").Value = -2764.5
$ws.Range("H67").Value = 1011.5
$ws.Range("I67").Value = 1011.5
$ws.Range("K67").Value = 3034.5
$ws.Range("M67").Value = -2098.5
$ws.Range("H75").Value = 799.5
$ws.Range("I75").Value = 799.5
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 2398.5
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -1400.5
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 799.5
$ws.Range("I78").Value = 799.5
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 7195.5
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -2203.5
$ws.Range("N78").ClearContents()
$ws.Range("H134").Value = 2844.0435
$ws.Range("I134").Value = 2989.1177
$ws.Range("K134").Value = 8967.3531
$ws.Range("M134").Value = -3897.3531

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 2201.2
$ws.Range("I21").Value = 2201.2
$ws.Range("K21").Value = 2201.2
$ws.Range("M21").Value = -2028.2
$ws.Range("H30").Value = 2201.2
$ws.Range("I30").Value = 2201.2
$ws.Range("K30").Value = 2201.2
$ws.Range("M30").Value = -2096.2
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2481.2222
$ws.Range("I132").Value = 2481.2222
$ws.Range("K132").Value = 7443.6666
$ws.Range("M132").Value = -4913.6666
$ws.Range("H136").Value = 30305338
$ws.Range("I136").Value = 2358.7144
$ws.Range("K136").Value = 7076.1432
$ws.Range("M136").Value = -4526.1432

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1150.5667
$ws.Range("I132").Value = 1091.52
$ws.Range("J132").Value = 1445.8
$ws.Range("K132").Value = 3274.56
$ws.Range("L132").Value = 4337.4
$ws.Range("M132").Value = -744.5599999999999
$ws.Range("N132").Value = -9397.4
$ws.Range("H136").Value = 1867.44
$ws.Range("I136").Value = 1734.4
$ws.Range("J136").Value = 2399.6
$ws.Range("K136").Value = 5203.200000000001
$ws.Range("L136").Value = 7198.799999999999
$ws.Range("M136").Value = -2653.200000000001
$ws.Range("N136").Value = -12298.8
